$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.866.55"
$ws.Range("E2").Value = "  -0.13%  "

$ws.Range("D3").Value = "1.736.05"
$ws.Range("E3").Value = "  -0.34%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.0000"
$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "242.39"
$ws.Range("E5").Value = "  +4.81%  "

$ws.Range("E6").Value = "  +0.09%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5183"
$ws.Range("E7").Value = "  -1.52%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2739"
$ws.Range("E8").Value = "  -0.68%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06161"
$ws.Range("E9").Value = "  +0.18%  "

$ws.Range("D10").Value = "1.739.93"
$ws.Range("E10").Value = "  +0.01%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07188"
$ws.Range("E11").Value = "  +1.32%  "

$ws.Range("E12").Value = "  -1.73%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6423"
$ws.Range("E13").Value = "  -0.33%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.613"
$ws.Range("E14").Value = "  +1.80%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "77.26"
$ws.Range("E15").Value = "  -0.48%  "

$ws.Range("E16").Value = "  +0.08%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.9999"
$ws.Range("E17").Value = "  +0.05%  "

$ws.Range("D18").Value = "25.901.29"
$ws.Range("E18").Value = "  +0.05%  "

$ws.Range("E19").Value = "  +1.99%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000006795"
$ws.Range("E20").Value = "  +1.55%  "

$ws.Range("D21").Value = "1.962.21"
$ws.Range("E21").Value = "  -0.02%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.281"
$ws.Range("E22").Value = "  -0.41%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.637"
$ws.Range("E23").Value = "  -1.52%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.291"
$ws.Range("E24").Value = "  +2.46%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "136.93"
$ws.Range("E25").Value = "  -2.47%  "

$ws.Range("E26").Value = "  -0.61%  "

$ws.Range("E27").Value = "  +0.42%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.774"
$ws.Range("E28").Value = "  -1.07%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "105.21"
$ws.Range("E29").Value = "  +2.40%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.964"
$ws.Range("E30").Value = "  +6.25%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08246"
$ws.Range("E31").Value = "  -0.76%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.647"

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04683"
$ws.Range("E33").Value = "  +3.38%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.657"
$ws.Range("E34").Value = "  +1.59%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9907"
$ws.Range("E35").Value = "  +1.14%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6195"
$ws.Range("E36").Value = "  -0.61%  "

$ws.Range("E37").Value = "  +0.30%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01601"
$ws.Range("E38").Value = "  +0.50%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.923"
$ws.Range("E39").Value = "  -1.30%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9997"
$ws.Range("E40").Value = "  +0.08%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "99.92"
$ws.Range("E41").Value = "  -0.23%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.3853"
$ws.Range("E42").Value = "  -0.62%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.7471"
$ws.Range("E43").Value = "  +1.95%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.008"
$ws.Range("E44").Value = "  -0.19%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1128"
$ws.Range("E45").Value = "  -0.08%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.252"
$ws.Range("E46").Value = "  -0.18%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "55.03"
$ws.Range("E47").Value = "  +2.49%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05214"
$ws.Range("E48").Value = "  -2.42%  "

$ws.Range("E49").Value = "  +1.46%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.565"
$ws.Range("E50").Value = "  -1.16%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3414"
$ws.Range("E51").Value = "  -0.54%  "
